# add noviat 70 modules
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OE Suppliers EUR")

# --- Update supplier / reference data (ASUSTeK -> OpenERP, PI14/xxxx -> PI12/xxxx) ---
$ws.Range("B3").Value = "OpenERP"
$ws.Range("A3").Value = "PI12/0101"
$ws.Range("A4").Value = "PI12/0250"
$ws.Range("D3").Value = 39824
$ws.Range("D4").Value = 39844
$ws.Range("H3").Value = "Open Invoice from supplier OpenERP"
$ws.Range("H4").Value = "Open Invoice from customer with reference '032/0029'"

# --- Move the active selection on the frozen pane from E8 to K7 ---
$ws.Range("K7").Select()
